$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A14").Value = "BanCoVan"
$ws.Range("B14").Value = 12345
$ws.Range("C14").Value = "super"
$ws.Range("C15").Select()
